$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J, matching the style of the existing
# header cells (bordered/bold/centered style used by B1:H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-8 for column I (all 1 except row 8 which is 4)
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 4

# Data rows 2-8 for column J
$ws.Range("J2").Value = 6
$ws.Range("J3").Value = 5
$ws.Range("J4").Value = 8
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 6
$ws.Range("J8").Value = 4
